$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("75:75").Insert()

$ws.Range("A75").Value = 4
$ws.Range("B75").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C75").Value = "Los Lagos"
$ws.Range("D75").Value = 44438
$ws.Range("E75").Value = 10
$ws.Range("F75").Value = 100112032
$ws.Range("G75").Value = "Zapallo italiano"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 70
$ws.Range("K75").Value = 16000
$ws.Range("L75").Value = 16000
$ws.Range("M75").Value = 16000
$ws.Range("N75").Value = "$/caja 50 unidades"
$ws.Range("O75").Value = "Región de Arica y Parinacota"
$ws.Range("P75").Value = 320
$ws.Range("Q75").Value = 50
$ws.Range("R75").Value = "Hortaliza"
